# Scheduled market-data refresh for Marilith_Profits workbook.
# Updates currentAveragePrice*/LevePrice*/LeveProfit* columns (H-N) across
# the ALC, ARM, BSM, CUL, GSM, LTW and WVR sheets with freshly pulled values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Cells.Item(64, 8).Value = 3750
$ws.Cells.Item(64, 10).Value = 3500
$ws.Cells.Item(64, 12).Value = 3500
$ws.Cells.Item(64, 14).Value = -3996

# Row 67
$ws.Cells.Item(67, 8).Value = 3750
$ws.Cells.Item(67, 10).Value = 3500
$ws.Cells.Item(67, 12).Value = 3500
$ws.Cells.Item(67, 14).Value = -5216

# Row 86
$ws.Cells.Item(86, 8).Value = 5237.8
$ws.Cells.Item(86, 9).Value = 2799
$ws.Cells.Item(86, 10).Value = 5847.5
$ws.Cells.Item(86, 11).Value = 2799
$ws.Cells.Item(86, 12).Value = 5847.5
$ws.Cells.Item(86, 13).Value = -1676
$ws.Cells.Item(86, 14).Value = -8093.5

# Row 89
$ws.Cells.Item(89, 8).Value = 5237.8
$ws.Cells.Item(89, 9).Value = 2799
$ws.Cells.Item(89, 10).Value = 5847.5
$ws.Cells.Item(89, 11).Value = 13995
$ws.Cells.Item(89, 12).Value = 29237.5
$ws.Cells.Item(89, 13).Value = -8379
$ws.Cells.Item(89, 14).Value = -40469.5

# Row 141
$ws.Cells.Item(141, 8).Value = 5957
$ws.Cells.Item(141, 9).Value = 5957
$ws.Cells.Item(141, 11).Value = 17871
$ws.Cells.Item(141, 13).Value = -12691

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5348.3213
$ws.Cells.Item(32, 9).Value = 4536.4727
$ws.Cells.Item(32, 10).Value = 50000
$ws.Cells.Item(32, 11).Value = 4536.4727
$ws.Cells.Item(32, 12).Value = 50000
$ws.Cells.Item(32, 13).Value = -4249.4727
$ws.Cells.Item(32, 14).Value = -50574

# Row 46
$ws.Cells.Item(46, 8).Value = 5268
$ws.Cells.Item(46, 9).Value = 5527
$ws.Cells.Item(46, 11).Value = 5527
$ws.Cells.Item(46, 13).Value = -5208

# Row 63
$ws.Cells.Item(63, 8).Value = 2582
$ws.Cells.Item(63, 9).Value = 1873
$ws.Cells.Item(63, 10).Value = 4000
$ws.Cells.Item(63, 11).Value = 1873
$ws.Cells.Item(63, 12).Value = 4000
$ws.Cells.Item(63, 13).Value = -1187
$ws.Cells.Item(63, 14).Value = -5372

# Row 66
$ws.Cells.Item(66, 8).Value = 2582
$ws.Cells.Item(66, 9).Value = 1873
$ws.Cells.Item(66, 10).Value = 4000
$ws.Cells.Item(66, 11).Value = 9365
$ws.Cells.Item(66, 12).Value = 20000
$ws.Cells.Item(66, 13).Value = -5933
$ws.Cells.Item(66, 14).Value = -26864

# Row 74
$ws.Cells.Item(74, 8).Value = 2403
$ws.Cells.Item(74, 10).Value = 2104.6667
$ws.Cells.Item(74, 12).Value = 2104.6667
$ws.Cells.Item(74, 14).Value = -3852.6667

# Row 77
$ws.Cells.Item(77, 8).Value = 2403
$ws.Cells.Item(77, 10).Value = 2104.6667
$ws.Cells.Item(77, 12).Value = 10523.3335
$ws.Cells.Item(77, 14).Value = -19259.3335

# Row 97
$ws.Cells.Item(97, 8).Value = 901.1667
$ws.Cells.Item(97, 9).Value = 852
$ws.Cells.Item(97, 10).Value = 999.5
$ws.Cells.Item(97, 11).Value = 852
$ws.Cells.Item(97, 12).Value = 999.5
$ws.Cells.Item(97, 13).Value = -356
$ws.Cells.Item(97, 14).Value = -1991.5

# Row 102
$ws.Cells.Item(102, 8).Value = 2499.25
$ws.Cells.Item(102, 9).Value = 2499.25
$ws.Cells.Item(102, 11).Value = 2499.25
$ws.Cells.Item(102, 13).Value = -877.25

# Row 122
$ws.Cells.Item(122, 8).Value = 2248.389
$ws.Cells.Item(122, 9).Value = 2298.353
$ws.Cells.Item(122, 10).Value = 1399
$ws.Cells.Item(122, 11).Value = 6895.059
$ws.Cells.Item(122, 12).Value = 4197
$ws.Cells.Item(122, 13).Value = -4445.059
$ws.Cells.Item(122, 14).Value = -9097

# Row 132
$ws.Cells.Item(132, 8).Value = 5539.364
$ws.Cells.Item(132, 9).Value = 5693.9
$ws.Cells.Item(132, 10).Value = 3994
$ws.Cells.Item(132, 11).Value = 17081.7
$ws.Cells.Item(132, 12).Value = 11982
$ws.Cells.Item(132, 13).Value = -14551.7
$ws.Cells.Item(132, 14).Value = -17042

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Cells.Item(82, 8).Value = 30286.5
$ws.Cells.Item(82, 10).Value = 39997.727
$ws.Cells.Item(82, 12).Value = 39997.727
$ws.Cells.Item(82, 14).Value = -40763.727

# Row 85
$ws.Cells.Item(85, 8).Value = 30286.5
$ws.Cells.Item(85, 10).Value = 39997.727
$ws.Cells.Item(85, 12).Value = 39997.727
$ws.Cells.Item(85, 14).Value = -42649.727

# Row 86
$ws.Cells.Item(86, 8).Value = 2580.2415
$ws.Cells.Item(86, 9).Value = 2512.4092
$ws.Cells.Item(86, 11).Value = 2512.4092
$ws.Cells.Item(86, 13).Value = -1389.4092

# Row 89
$ws.Cells.Item(89, 8).Value = 2580.2415
$ws.Cells.Item(89, 9).Value = 2512.4092
$ws.Cells.Item(89, 11).Value = 12562.046
$ws.Cells.Item(89, 13).Value = -6946.046

# Row 94
$ws.Cells.Item(94, 8).Value = 1533.5454
$ws.Cells.Item(94, 9).Value = 1596.4736
$ws.Cells.Item(94, 11).Value = 1596.4736
$ws.Cells.Item(94, 13).Value = -1145.4736

# Row 99
$ws.Cells.Item(99, 8).Value = 3145.6155
$ws.Cells.Item(99, 9).Value = 3324.5
$ws.Cells.Item(99, 11).Value = 3324.5
$ws.Cells.Item(99, 13).Value = -1826.5

# Row 107
$ws.Cells.Item(107, 8).Value = 1713.7142
$ws.Cells.Item(107, 9).Value = 1166
$ws.Cells.Item(107, 11).Value = 1166
$ws.Cells.Item(107, 13).Value = 754

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Cells.Item(14, 8).Value = 29899.924
$ws.Cells.Item(14, 9).Value = 29899.924
$ws.Cells.Item(14, 11).Value = 89699.772
$ws.Cells.Item(14, 13).Value = -89526.772

# Row 50
$ws.Cells.Item(50, 8).Value = 515
$ws.Cells.Item(50, 9).Value = 515
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 1545
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = -1064
$ws.Cells.Item(50, 14).ClearContents()

# Row 53
$ws.Cells.Item(53, 8).Value = 515
$ws.Cells.Item(53, 9).Value = 515
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 1545
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = -1064
$ws.Cells.Item(53, 14).ClearContents()

# Row 80
$ws.Cells.Item(80, 8).Value = 1330
$ws.Cells.Item(80, 9).Value = 1397
$ws.Cells.Item(80, 11).Value = 4191
$ws.Cells.Item(80, 13).Value = -3255

# Row 83
$ws.Cells.Item(83, 8).Value = 1330
$ws.Cells.Item(83, 9).Value = 1397
$ws.Cells.Item(83, 11).Value = 12573
$ws.Cells.Item(83, 13).Value = -7893

# Row 92
$ws.Cells.Item(92, 8).Value = 798.6667
$ws.Cells.Item(92, 10).Value = 698
$ws.Cells.Item(92, 12).Value = 2094
$ws.Cells.Item(92, 14).Value = -4590

# Row 107
$ws.Cells.Item(107, 8).Value = 1641.7142
$ws.Cells.Item(107, 9).Value = 1641.7142
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 4925.142599999999
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -3005.142599999999
$ws.Cells.Item(107, 14).ClearContents()

# Row 113
$ws.Cells.Item(113, 8).Value = 1109.1
$ws.Cells.Item(113, 10).Value = 1176.7778
$ws.Cells.Item(113, 12).Value = 3530.3334
$ws.Cells.Item(113, 14).Value = -7870.3334

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 6287031.5
$ws.Cells.Item(122, 9).Value = 7388272
$ws.Cells.Item(122, 11).Value = 22164816
$ws.Cells.Item(122, 13).Value = -22162366

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 4077.3333
$ws.Cells.Item(7, 10).Value = 3399.6667
$ws.Cells.Item(7, 12).Value = 3399.6667
$ws.Cells.Item(7, 14).Value = -3623.6667

# Row 13
$ws.Cells.Item(13, 8).Value = 6666
$ws.Cells.Item(13, 10).Value = 6666
$ws.Cells.Item(13, 12).Value = 6666
$ws.Cells.Item(13, 14).Value = -6946

# Row 22
$ws.Cells.Item(22, 8).Value = 1440.409
$ws.Cells.Item(22, 9).Value = 1121.7222
$ws.Cells.Item(22, 10).Value = 2874.5
$ws.Cells.Item(22, 11).Value = 1121.7222
$ws.Cells.Item(22, 12).Value = 2874.5
$ws.Cells.Item(22, 13).Value = -826.7221999999999
$ws.Cells.Item(22, 14).Value = -3464.5

# Row 27
$ws.Cells.Item(27, 8).Value = 1440.409
$ws.Cells.Item(27, 9).Value = 1121.7222
$ws.Cells.Item(27, 10).Value = 2874.5
$ws.Cells.Item(27, 11).Value = 1121.7222
$ws.Cells.Item(27, 12).Value = 2874.5
$ws.Cells.Item(27, 13).Value = -1014.7222
$ws.Cells.Item(27, 14).Value = -3088.5

# Row 40
$ws.Cells.Item(40, 8).Value = 7400
$ws.Cells.Item(40, 9).Value = 6100
$ws.Cells.Item(40, 10).Value = 10000
$ws.Cells.Item(40, 11).Value = 6100
$ws.Cells.Item(40, 12).Value = 10000
$ws.Cells.Item(40, 13).Value = -5964
$ws.Cells.Item(40, 14).Value = -10272

# Row 46
$ws.Cells.Item(46, 8).Value = 4499.1665
$ws.Cells.Item(46, 9).Value = 4750
$ws.Cells.Item(46, 10).Value = 3997.5
$ws.Cells.Item(46, 11).Value = 4750
$ws.Cells.Item(46, 12).Value = 3997.5
$ws.Cells.Item(46, 13).Value = -4562
$ws.Cells.Item(46, 14).Value = -4373.5

# Row 61
$ws.Cells.Item(61, 8).Value = 6826.909
$ws.Cells.Item(61, 9).Value = 5789
$ws.Cells.Item(61, 11).Value = 5789
$ws.Cells.Item(61, 13).Value = -5587

# Row 68
$ws.Cells.Item(68, 8).Value = 4225
$ws.Cells.Item(68, 9).Value = 3950
$ws.Cells.Item(68, 11).Value = 3950
$ws.Cells.Item(68, 13).Value = -3201

# Row 71
$ws.Cells.Item(71, 8).Value = 4225
$ws.Cells.Item(71, 9).Value = 3950
$ws.Cells.Item(71, 11).Value = 19750
$ws.Cells.Item(71, 13).Value = -16006

# Row 113
$ws.Cells.Item(113, 8).Value = 6826.909
$ws.Cells.Item(113, 9).Value = 5789
$ws.Cells.Item(113, 11).Value = 5789
$ws.Cells.Item(113, 13).Value = -3619

# Row 126
$ws.Cells.Item(126, 8).Value = 4077.3333
$ws.Cells.Item(126, 10).Value = 3399.6667
$ws.Cells.Item(126, 12).Value = 10199.0001
$ws.Cells.Item(126, 14).Value = -15139.0001

# Row 132
$ws.Cells.Item(132, 8).Value = 20039.7
$ws.Cells.Item(132, 9).Value = 20877.445
$ws.Cells.Item(132, 10).Value = 12500
$ws.Cells.Item(132, 11).Value = 62632.335
$ws.Cells.Item(132, 12).Value = 37500
$ws.Cells.Item(132, 13).Value = -60102.335
$ws.Cells.Item(132, 14).Value = -42560

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 674.5
$ws.Cells.Item(107, 9).Value = 500
$ws.Cells.Item(107, 10).Value = 732.6667
$ws.Cells.Item(107, 11).Value = 1500
$ws.Cells.Item(107, 12).Value = 2198.0001
$ws.Cells.Item(107, 13).Value = 420
$ws.Cells.Item(107, 14).Value = -6038.0001

# Row 132
$ws.Cells.Item(132, 8).Value = 1092.625
$ws.Cells.Item(132, 9).Value = 983.0833
$ws.Cells.Item(132, 11).Value = 2949.2499
$ws.Cells.Item(132, 13).Value = -419.2498999999998

